# Fixed ND and BN processing
# - Course codes on "N-SECOND-YEAR-FIRST-SEMESTER" lose the stray space
#   ("GNS 211" -> "GNS211", etc.)
# - That sheet becomes the active tab/selection, replacing
#   "N-FIRST-YEAR-SECOND-SEMESTER" as the selected one.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("N-SECOND-YEAR-FIRST-SEMESTER")

$ws.Range("A2").Value = "GNS211"
$ws.Range("A3").Value = "GNS212"
$ws.Range("A4").Value = "GNS213"
$ws.Range("A5").Value = "GNS214"
$ws.Range("A6").Value = "GNS215"
$ws.Range("A7").Value = "GNS216"
$ws.Range("A8").Value = "GNS217"
$ws.Range("A9").Value = "GNS218"

# Make this sheet the active/selected tab, with its own new active cell.
$ws.Activate()
$ws.Range("B15").Select()
